# Add two new completed books to the "Completed" sheet:
#   - The Years That Matter Most, by Paul Tough
#   - Springfield Confidential, by Mike Reiss

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Row 83: The Years That Matter Most
$ws.Range("A83").Value = "The Years That Matter Most"
$ws.Range("B83").Value = "Paul Tough"
$ws.Range("C83").Value2 = 43984
$ws.Range("D83").Value2 = 43986
$ws.Range("E83").Value = "college;SAT;admissions;discrimination;equality"
$ws.Range("F83").Value = "Audio"
$ws.Range("G83").Value = "12 Hours 39 Mins"

# Row 84: Springfield Confidential
$ws.Range("A84").Value = "Springfield Confidential"
$ws.Range("B84").Value = "Mike Reiss"
$ws.Range("C84").Value2 = 43986
$ws.Range("D84").Value2 = 43987
$ws.Range("E84").Value = "simpsons;tv;writing;comedy;sitcoms"
$ws.Range("F84").Value = "Audio"
$ws.Range("G84").Value = "7 Hours 34 Mins"

# Copy the existing date-cell formatting down onto the new rows
$ws.Range("C82:D82").Copy()
$ws.Range("C83:D84").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A85").Select()
